# Auto-generated Word COM-interop edit script
$d = $word.ActiveDocument

function ReplaceInRange($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function AppendAtParagraphEnd($para, $text) {
    $endPos = $para.Range.End
    $insertPoint = $d.Range($endPos - 1, $endPos - 1)
    $insertPoint.InsertAfter($text)
}

# --- Title ---
$p1 = $d.Paragraphs.Item(1)
ReplaceInRange $p1.Range 'The Unseen Enemy: Battling Antimicrobial Resistance' 'Understanding the Past: History and Its Significance for High School Students'

# --- Author name ---
$p2 = $d.Paragraphs.Item(2)
ReplaceInRange $p2.Range 'Alexciaa Mendez' 'Rachel Williams'

# --- Email ---
$p3 = $d.Paragraphs.Item(3)
ReplaceInRange $p3.Range 'alexciaamendez@gmail' 'rachel'
ReplaceInRange $p3.Range 'com' 'williams@highschools'
AppendAtParagraphEnd $p3 '.'
AppendAtParagraphEnd $p3 'edu'

# --- Body paragraph sentence replacements ---
$p5 = $d.Paragraphs.Item(5)
ReplaceInRange $p5.Range 'In the realm of healthcare, there lies a daunting threat that has the potential to undermine modern medicine''s capabilities: antimicrobial resistance (AMR)' 'Perched on the threshold of knowledge, High School Students are presented with a panorama of disciplines begging to be explored'
ReplaceInRange $p5.Range ' This phenomenon arises when microorganisms, such as bacteria, develop the ability to evade the effects of antimicrobial agents, rendering treatments ineffective' ' History, a subject of paramount significance, beckons with its siren call, whispering tales of bygone eras and the shadows they cast upon our present'
ReplaceInRange $p5.Range ' The consequences are dire and far-reaching, jeopardizing our capacity to combat infections and threatening the very foundation of healthcare' ' This introduction delves into the fascinating realm of history, elucidating its intrinsic value and unraveling its relevance to high school students'' lives'
ReplaceInRange $p5.Range 'AMR is not a futuristic concept; it is a pressing global health crisis already impacting countless lives worldwide' 'History is the compass by which we navigate the complexities of the present'
ReplaceInRange $p5.Range ' Each year, millions of people succumb to infections caused by drug-resistant microbes, and the numbers continue to rise alarmingly' ' By comprehending the past, we gain profound insights into the myriad cultural, political, and societal forces that have shaped our world'
ReplaceInRange $p5.Range ' The emergence of "superbugs," organisms resistant to multiple antibiotics, poses a significant challenge, limiting treatment options and increasing the likelihood of prolonged illnesses, severe complications, and even death' ' The echoes of ancient civilizations, the reverberations of past conflicts, and the whispers of forgotten triumphs weave together an intricate tapestry of understanding, empowering us to make informed decisions and navigate the ever-changing landscape of life'
ReplaceInRange $p5.Range 'The widespread use and misuse of antibiotics have contributed to the acceleration of AMR' 'Furthermore, history is a testament to human resilience and ingenuity'
ReplaceInRange $p5.Range ' The excessive and inappropriate prescription of antibiotics, often for viral infections where they are ineffective, has created a breeding ground for resistant bacteria' ' It unveils narratives of perseverance, innovation, and indomitable spirit, inspiring us to surmount challenges and reach for greatness'
ReplaceInRange $p5.Range ' Furthermore, the overuse of antibiotics in agriculture, to promote growth in livestock and prevent disease, has exacerbated the problem' ' From the architectural wonders of ancient civilizations to the scientific breakthroughs of modern times, history reminds us of our potential to create a better world and to shape the course of destiny'
ReplaceInRange $p5.Range ' The imprudent use of antimicrobials has fueled AMR, driving the evolution of resistant microorganisms that can spread within and between human, animal, and environmental populations' ' It is in these stories of courage, resilience, and vision that we find the impetus to make a difference in our own lives and the lives of others'

# --- Append new trailing content to body paragraph ---
AppendAtParagraphEnd $p5 "`v`vBody:`v`vHistory, as an academic discipline, offers an unparalleled opportunity to develop critical thinking and problem-solving skills. By examining primary and secondary sources, students learn to evaluate evidence, construct compelling arguments, and engage in meaningful debates. These skills are not only essential for academic success but also for navigating the complexities of modern life. History teaches us to question assumptions, think critically, and make informed decisions based on evidence - skills that are invaluable in the workplace and in civic life.`v`vIn addition to its practical applications, history also cultivates a profound sense of empathy and cultural awareness. By exploring the lives and experiences of people from different times and places, students develop a deeper understanding of the human condition. They learn to appreciate the richness and diversity of human cultures, to respect different perspectives, and to recognize the common threads that bind humanity together. This empathy and cultural awareness are essential for creating a more just, equitable, and peaceful world.`v`vMoreover, history provides a profound sense of identity and belonging. It helps us understand our roots, our heritage, and the traditions that have shaped us. By studying history, we gain a deeper appreciation for our culture and our place in the world. It fosters a sense of pride and patriotism while encouraging us to critically examine the past and work towards a better future. This sense of identity and belonging is vital for our emotional well-being and for our ability to contribute positively to society."

# --- Summary paragraph ---
$p7 = $d.Paragraphs.Item(7)
ReplaceInRange $p7.Range 'Antimicrobial resistance, a formidable adversary challenging modern medicine, has become a global health crisis' 'In conclusion, history is an indispensable subject for High School Students'
ReplaceInRange $p7.Range ' The emergence of drug-resistant microorganisms threatens the efficacy of antibiotics and poses significant risks to public health' ' It offers a profound understanding of the past, cultivates critical thinking and problem-solving skills, promotes empathy and cultural awareness, and fosters a sense of identity and belonging'
ReplaceInRange $p7.Range ' Urgent action is required to address this pressing issue through responsible antibiotic use, enhanced infection prevention and control measures, and the development of novel antimicrobial agents' ' By studying history, students gain the knowledge and skills necessary for success in college, career, and life'
ReplaceInRange $p7.Range ' By working collectively, healthcare professionals, policymakers, industry leaders, and the public can combat AMR and safeguard the future of effective healthcare' ' They become informed citizens, capable of making informed decisions and contributing positively to society'

# Append two new trailing runs (period + new final sentence) before the last unchanged period
AppendAtParagraphEnd $p7 '. As they navigate the challenges and opportunities of the 21st century, history will serve as their compass, their inspiration, and their guide'

# --- Add a trailing empty paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Output "done"
